$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.680.36"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "3.520.63"
$ws.Range("E3").Value = "  -2.52%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'619.12"
$ws.Range("E5").Value = "  +2.88%  "
$ws.Range("D6").Value = "'172.98"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "  -1.28%  "
$ws.Range("D8").Value = "3.517.60"
$ws.Range("E8").Value = "  -2.46%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'0.199"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").Value = "'7.13"
$ws.Range("E11").Value = "  -4.39%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "'46.40"
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "4.088.78"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").Value = "'8.36"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "'609.32"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").Value = "3.518.86"
$ws.Range("E18").Value = "  -2.76%  "
$ws.Range("D19").Value = "70.727.08"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").Value = "'0.883"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "'9.19"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("E24").Value = "  -3.43%  "
$ws.Range("D25").Value = "'97.77"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  -2.35%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -4.31%  "
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").Value = "'9.08"
$ws.Range("E30").Value = "  -2.89%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'8.08"
$ws.Range("E31").Value = "  -5.82%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").Value = "'2.99"
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("D34").Value = "'641.28"
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("E35").Value = "  -6.59%  "
$ws.Range("E36").Value = "  -2.93%  "
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").Value = "'0.0485"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "'3.41"
$ws.Range("E39").Value = "  -9.34%  "
$ws.Range("D40").Value = "'56.82"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "3.351.71"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("D45").Value = "'2.97"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("E46").Value = "  -5.01%  "
$ws.Range("D47").Value = "'31.81"
$ws.Range("E47").Value = "  -4.59%  "
$ws.Range("E48").Value = "  -6.67%  "
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").Value = "'134.12"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("E51").Value = "  -0.01%  "